$wb = $excel.ActiveWorkbook

# The two new sheets are clones of the last existing sheet (Bus_Makhulu_r),
# renamed for a pair of new multi-axle truck variants, each with its own
# distinct label cell (H3) feeding the sheet's unique shared-string entry.
$src = $wb.Worksheets.Item("Bus_Makhulu_r")

# --- New sheet 1: Truck_Amandla_A2 -----------------------------------
$src.Copy([System.Reflection.Missing]::Value, $src)
$ws7 = $wb.Worksheets.Item($src.Index + 1)
$ws7.Name = "Truck_Amandla_A2"
$ws7.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A2"
$ws7.Activate() | Out-Null
$ws7.Range("J9").Select() | Out-Null

# --- New sheet 2: Truck_Amandla_A3 (cloned after the first new sheet) -
$ws7.Copy([System.Reflection.Missing]::Value, $ws7)
$ws8 = $wb.Worksheets.Item($ws7.Index + 1)
$ws8.Name = "Truck_Amandla_A3"
$ws8.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A3"

# The newly appended last sheet ends up the active/selected tab.
$ws8.Activate() | Out-Null
$ws8.Range("E13").Select() | Out-Null
